$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Delete EBO row (original row 22). EBO will be merged into EPB. ---
$ws.Rows.Item(22).Delete()

# --- Step 2: Row 29 now holds EMG (shifted up from 30 by the EBO deletion). ---
# Rename it to EMR and set its values to the merged sums of ENF + EMG.
$ws.Cells.Item(29, 1).Value = "EMR"
$ws.Cells.Item(29, 2).Value = 120214
$ws.Cells.Item(29, 3).Value = 121806
$ws.Cells.Item(29, 4).Value = 123205
$ws.Cells.Item(29, 5).Value = 124820
$ws.Cells.Item(29, 6).Value = 126183
$ws.Cells.Item(29, 7).Value = 127403
$ws.Cells.Item(29, 8).Value = 127348
$ws.Cells.Item(29, 9).Value = 127777
$ws.Cells.Item(29, 10).Value = 12722
$ws.Cells.Item(29, 11).Value = 12686

# --- Step 3: Delete the ENF row (now at row 35 after the EBO deletion shift). ---
# Its data has already been folded into the EMR row above.
$ws.Rows.Item(35).Delete()

# --- Step 4: Row 35 now holds EPB (shifted up after the ENF deletion). ---
# Set its values to the merged sums of EBO + EPB.
$ws.Cells.Item(35, 2).Value = 257841
$ws.Cells.Item(35, 3).Value = 266030
$ws.Cells.Item(35, 4).Value = 271581
$ws.Cells.Item(35, 5).Value = 274088
$ws.Cells.Item(35, 6).Value = 277374
$ws.Cells.Item(35, 7).Value = 280411
$ws.Cells.Item(35, 8).Value = 283615
$ws.Cells.Item(35, 9).Value = 286986
$ws.Cells.Item(35, 10).Value = 287200
$ws.Cells.Item(35, 11).Value = 263839

# --- Step 5: Add new column L (2023 data) header, formatted like the other year headers. ---
$ws.Range("L1").NumberFormat = "@"
$ws.Range("L1").Value = "2023"
$ws.Range("L1").NumberFormat = "General"

# --- Step 6: Add the 2023 values for every distributor row in the final layout (rows 2-53). ---
$ws.Cells.Item(2, 12).Value = 130067
$ws.Cells.Item(3, 12).Value = 20996
$ws.Cells.Item(4, 12).Value = 131058
$ws.Cells.Item(5, 12).Value = 270962
$ws.Cells.Item(6, 12).Value = 700737
$ws.Cells.Item(7, 12).Value = 385615
$ws.Cells.Item(8, 12).Value = 1475318
$ws.Cells.Item(9, 12).Value = 8541
$ws.Cells.Item(10, 12).Value = 5699
$ws.Cells.Item(11, 12).Value = 707785
$ws.Cells.Item(12, 12).Value = 5386
$ws.Cells.Item(13, 12).Value = 872385
$ws.Cells.Item(14, 12).Value = 187258
$ws.Cells.Item(15, 12).Value = 369433
$ws.Cells.Item(16, 12).Value = 107394
$ws.Cells.Item(17, 12).Value = 53213
$ws.Cells.Item(18, 12).Value = 11065
$ws.Cells.Item(19, 12).Value = 3748
$ws.Cells.Item(20, 12).Value = 7455
$ws.Cells.Item(21, 12).Value = 45242
$ws.Cells.Item(22, 12).Value = 340267
$ws.Cells.Item(23, 12).Value = 176886
$ws.Cells.Item(24, 12).Value = 548
$ws.Cells.Item(25, 12).Value = 1850
$ws.Cells.Item(26, 12).Value = 379649
$ws.Cells.Item(27, 12).Value = 8132
$ws.Cells.Item(28, 12).Value = 47022
$ws.Cells.Item(29, 12).Value = 12656
$ws.Cells.Item(30, 12).Value = 178852
$ws.Cells.Item(31, 12).Value = 306133
$ws.Cells.Item(32, 12).Value = 646116
$ws.Cells.Item(33, 12).Value = 218938
$ws.Cells.Item(34, 12).Value = 462124
$ws.Cells.Item(35, 12).Value = 258417
$ws.Cells.Item(36, 12).Value = 115317
$ws.Cells.Item(37, 12).Value = 452126
$ws.Cells.Item(38, 12).Value = 219504
$ws.Cells.Item(39, 12).Value = 370576
$ws.Cells.Item(40, 12).Value = 228574
$ws.Cells.Item(41, 12).Value = 171101
$ws.Cells.Item(42, 12).Value = 77373
$ws.Cells.Item(43, 12).Value = 121150
$ws.Cells.Item(44, 12).Value = 89034
$ws.Cells.Item(45, 12).Value = 1893
$ws.Cells.Item(46, 12).Value = 297023
$ws.Cells.Item(47, 12).Value = 1970
$ws.Cells.Item(48, 12).Value = 245763
$ws.Cells.Item(49, 12).Value = 1852
$ws.Cells.Item(50, 12).Value = 452728
$ws.Cells.Item(51, 12).Value = 28712
$ws.Cells.Item(52, 12).Value = 11695
$ws.Cells.Item(53, 12).Value = 6331
